# TCA.xlsx update — add a new "Scenario4" column and rescale the scenario
# shock sizes. Net effect (vs. the original B:E layout) is equivalent to
# inserting a new column before B, filling it with a fresh set of shock
# values, touching up a handful of the old B/C columns (now C/D), and
# appending a new F column that repeats the old closing "*1.5" formula one
# column further to the right.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for the appended scenario column.
$ws.Range("F1").Value = "Scenario4"

# --- Column B (new scenario shocks) ---------------------------------
$ws.Range("B2:B9").Value = 0.00002
$ws.Range("B10:B21").Value = 0.00005

# --- Column C ---------------------------------------------------------
$ws.Range("C2:C3").Value = 0.00002
$ws.Range("C4:C7").Value = 0.00005
$ws.Range("C8:C11").Value = 0.00007
$ws.Range("C12:C21").Value = 0.0001

# --- Column D -----------------------------------------------------------
$ws.Range("D2:D9").Value = 0.00005
$ws.Range("D10:D21").Value = 0.0001

# --- Column E (old column D content, shifted one column right) ----------
$ws.Range("E2:E8").Value = 0.0001
$ws.Range("E9").Formula = "=0.0001*5"
$ws.Range("E10:E11").Formula = "=0.0001*5"
$ws.Range("E12").Formula = "=0.0001*10"
$ws.Range("E13:E21").Formula = "=0.0001*10"

# --- Column F (old closing formula, shifted one column right, now off E) --
$ws.Range("F2").Formula = "=E2*1.5"
$ws.Range("F3:F21").Formula = "=E3*1.5"

# Recalculate so every cell carries a fresh cached value.
$wb.Application.Calculate()

# Restore the active selection to where the user last clicked.
$ws.Range("G10").Select()
